$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number ("179.30", "0.625", ...)
# must be forced to Text format first, otherwise Excel auto-converts them to
# a numeric value (dropping the trailing zero / thousands dots) when assigned
# via .Value. ClearFormats() afterwards drops the temporary "@" style again so
# the cell keeps the workbook default style (matches the source, which never set
# an explicit style on these cells).
$numericLooking = @(
    "D5", "D6", "D7", "D11", "D15", "D18", "D19", "D20", "D22", "D28", "D33", "D36", "D38", "D41", "D43", "D44", "D46", "D47", "D48", "D49", "D51"
)
foreach ($addr in $numericLooking) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '66.112.73'
$ws.Range("E2").Value = '  -1.85%  '
$ws.Range("D3").Value = '3.278.66'
$ws.Range("E3").Value = '  -1.28%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '575.62'
$ws.Range("E5").Value = '  -0.40%  '
$ws.Range("D6").Value = '179.30'
$ws.Range("E6").Value = '  -3.58%  '
$ws.Range("D7").Value = '0.625'
$ws.Range("E7").Value = '  +3.19%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  -3.14%  '
$ws.Range("E10").Value = '  +0.88%  '
$ws.Range("D11").Value = '0.400'
$ws.Range("E11").Value = '  -1.68%  '
$ws.Range("D12").Value = '3.848.48'
$ws.Range("E12").Value = '  -1.28%  '
$ws.Range("E13").Value = '  -3.85%  '
$ws.Range("D14").Value = '66.176.40'
$ws.Range("E14").Value = '  -2.03%  '
$ws.Range("D15").Value = '26.49'
$ws.Range("E15").Value = '  -3.25%  '
$ws.Range("E16").Value = '  -2.71%  '
$ws.Range("D17").Value = '3.280.06'
$ws.Range("E17").Value = '  -1.07%  '
$ws.Range("D18").Value = '431.79'
$ws.Range("E18").Value = '  -2.61%  '
$ws.Range("D19").Value = '5.56'
$ws.Range("E19").Value = '  -1.94%  '
$ws.Range("D20").Value = '13.12'
$ws.Range("E20").Value = '  -3.31%  '
$ws.Range("E21").Value = '  -4.13%  '
$ws.Range("D22").Value = '71.93'
$ws.Range("E22").Value = '  -3.04%  '
$ws.Range("E23").Value = '  +0.16%  '
$ws.Range("D24").Value = '3.414.52'
$ws.Range("E24").Value = '  -1.47%  '
$ws.Range("E25").Value = '  -1.49%  '
$ws.Range("E26").Value = '  +3.33%  '
$ws.Range("E27").Value = '  -5.24%  '
$ws.Range("D28").Value = '8.88'
$ws.Range("E28").Value = '  -1.92%  '
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("E30").Value = '  -2.19%  '
$ws.Range("E31").Value = '  -2.81%  '
$ws.Range("D33").Value = '5.15'
$ws.Range("E33").Value = '  -3.40%  '
$ws.Range("E34").Value = '  -3.59%  '
$ws.Range("E35").Value = '  -4.60%  '
$ws.Range("D36").Value = '158.13'
$ws.Range("E36").Value = '  -2.57%  '
$ws.Range("E37").Value = '  -5.78%  '
$ws.Range("D38").Value = '26.59'
$ws.Range("E38").Value = '  -2.48%  '
$ws.Range("E39").Value = '  -3.12%  '
$ws.Range("D40").Value = '2.764.82'
$ws.Range("E40").Value = '  -0.84%  '
$ws.Range("D41").Value = '0.775'
$ws.Range("E41").Value = '  -1.77%  '
$ws.Range("E42").Value = '  -3.96%  '
$ws.Range("D43").Value = '40.26'
$ws.Range("E43").Value = '  -0.17%  '
$ws.Range("D44").Value = '6.03'
$ws.Range("E44").Value = '  -3.32%  '
$ws.Range("E45").Value = '  -2.36%  '
$ws.Range("B46").Value = 'dogwifhat'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D46").Value = '2.30'
$ws.Range("E46").Value = '  -3.97%  '
$ws.Range("B47").Value = 'Bittensor'
$ws.Range("C47").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D47").Value = '321.07'
$ws.Range("E47").Value = '  -1.87%  '
$ws.Range("D48").Value = '23.26'
$ws.Range("E48").Value = '  -5.84%  '
$ws.Range("D49").Value = '0.0266'
$ws.Range("E49").Value = '  -2.42%  '
$ws.Range("E50").Value = '  +1.87%  '
$ws.Range("D51").Value = '1.00'
$ws.Range("E51").Value = '  +0.00%  '

foreach ($addr in $numericLooking) {
    $ws.Range($addr).ClearFormats()
}
